$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data point (Fecha=45211) is inserted as the newest entry for
# this subset, ahead of the existing history. Insert two blank rows at the
# top of the data block (row 453) so every following row shifts down by two
# (the two rows that fall off the bottom of the previous range become the
# new last rows, rows 563-564).
$ws.Range("A453:A454").EntireRow.Insert()

# Row 453: Apio / Americana (o) / Primera, Fecha 45211
$ws.Cells.Item(453, 1).Value = 9
$ws.Cells.Item(453, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(453, 3).Value = "Metropolitana"
$ws.Cells.Item(453, 4).Value = 45211
$ws.Cells.Item(453, 5).Value = 13
$ws.Cells.Item(453, 6).Value = 100112017
$ws.Cells.Item(453, 7).Value = "Apio"
$ws.Cells.Item(453, 8).Value = "Americana (o)"
$ws.Cells.Item(453, 9).Value = "Primera"
$ws.Cells.Item(453, 10).Value = 70
$ws.Cells.Item(453, 11).Value = 6000
$ws.Cells.Item(453, 12).Value = 7000
$ws.Cells.Item(453, 13).Value = 6500
$ws.Cells.Item(453, 14).Value = "`$/docena de matas"
$ws.Cells.Item(453, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(453, 16).Value = 1083
$ws.Cells.Item(453, 17).Value = 6
$ws.Cells.Item(453, 18).Value = "Hortaliza"

# Row 454: Apio / Americana (o) / Segunda, Fecha 45211
$ws.Cells.Item(454, 1).Value = 9
$ws.Cells.Item(454, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(454, 3).Value = "Metropolitana"
$ws.Cells.Item(454, 4).Value = 45211
$ws.Cells.Item(454, 5).Value = 13
$ws.Cells.Item(454, 6).Value = 100112017
$ws.Cells.Item(454, 7).Value = "Apio"
$ws.Cells.Item(454, 8).Value = "Americana (o)"
$ws.Cells.Item(454, 9).Value = "Segunda"
$ws.Cells.Item(454, 10).Value = 52
$ws.Cells.Item(454, 11).Value = 5000
$ws.Cells.Item(454, 12).Value = 5000
$ws.Cells.Item(454, 13).Value = 5000
$ws.Cells.Item(454, 14).Value = "`$/docena de matas"
$ws.Cells.Item(454, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(454, 16).Value = 833
$ws.Cells.Item(454, 17).Value = 6
$ws.Cells.Item(454, 18).Value = "Hortaliza"
